# Auto-generated edit script
# Applies updated market-price / profit figures to the Aegis_Profits workbook
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the scheduled runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2
$ws.Range("H2").Value = 547.1
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 610.1429000000001
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 610.1429000000001
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -836.1429000000001

# ALC row 38
$ws.Range("H38").Value = 1897867
$ws.Range("I38").Value = 3225908.8
$ws.Range("J38").Value = 664.5714
$ws.Range("K38").Value = 9677726.399999999
$ws.Range("L38").Value = 1993.7142
$ws.Range("M38").Value = -9677354.399999999
$ws.Range("N38").Value = -2737.7142

# ALC row 51
$ws.Range("H51").Value = 14289.444
$ws.Range("I51").Value = 27100
$ws.Range("J51").Value = 4041
$ws.Range("K51").Value = 27100
$ws.Range("L51").Value = 4041
$ws.Range("M51").Value = -26616
$ws.Range("N51").Value = -5009

# ALC row 58
$ws.Range("H58").Value = 1264050.5
$ws.Range("I58").Value = 2525601
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 7576803
$ws.Range("L58").Value = 7500
$ws.Range("M58").Value = -7576653
$ws.Range("N58").Value = -7800

# ALC row 87
$ws.Range("H87").Value = 25547.875
$ws.Range("J87").Value = 25547.875
$ws.Range("L87").Value = 25547.875
$ws.Range("N87").Value = -28043.875

# ALC row 90
$ws.Range("H90").Value = 25547.875
$ws.Range("J90").Value = 25547.875
$ws.Range("L90").Value = 76643.625
$ws.Range("N90").Value = -89123.625

$ws = $wb.Worksheets.Item("ARM")
# ARM row 56
$ws.Range("H56").Value = 20000
$ws.Range("I56").Value = 20000
$ws.Range("K56").Value = 20000
$ws.Range("M56").Value = -19258

# ARM row 74
$ws.Range("H74").Value = 1459.9298
$ws.Range("I74").Value = 1008.44116
$ws.Range("J74").Value = 2127.348
$ws.Range("K74").Value = 1008.44116
$ws.Range("L74").Value = 2127.348
$ws.Range("M74").Value = -134.44116
$ws.Range("N74").Value = -3875.348

# ARM row 77
$ws.Range("H77").Value = 1459.9298
$ws.Range("I77").Value = 1008.44116
$ws.Range("J77").Value = 2127.348
$ws.Range("K77").Value = 5042.2058
$ws.Range("L77").Value = 10636.74
$ws.Range("M77").Value = -674.2057999999997
$ws.Range("N77").Value = -19372.74

# ARM row 98
$ws.Range("H98").Value = 18941
$ws.Range("J98").Value = 18941
$ws.Range("L98").Value = 18941
$ws.Range("N98").Value = -24931

$ws = $wb.Worksheets.Item("BSM")
# BSM row 54
$ws.Range("H54").Value = 6023.5
$ws.Range("I54").Value = 2410.25
$ws.Range("J54").Value = 13250
$ws.Range("K54").Value = 2410.25
$ws.Range("L54").Value = 13250
$ws.Range("M54").Value = -1926.25
$ws.Range("N54").Value = -14218

# BSM row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# BSM row 107
$ws.Range("H107").Value = 15172517
$ws.Range("I107").Value = 18543070
$ws.Range("J107").Value = 5029.5
$ws.Range("K107").Value = 18543070
$ws.Range("L107").Value = 5029.5
$ws.Range("M107").Value = -18541150
$ws.Range("N107").Value = -8869.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 683.2
$ws.Range("I22").Value = 366.4
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 366.4
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -16.39999999999998
$ws.Range("N22").Value = -1700

# CRP row 28
$ws.Range("H28").Value = 8999
$ws.Range("J28").Value = 8999
$ws.Range("L28").Value = 8999
$ws.Range("N28").Value = -9489

# CRP row 31
$ws.Range("H31").Value = 18813.83
$ws.Range("I31").Value = 46584.227
$ws.Range("J31").Value = 2301.7026
$ws.Range("K31").Value = 46584.227
$ws.Range("L31").Value = 2301.7026
$ws.Range("M31").Value = -46289.227
$ws.Range("N31").Value = -2891.7026

# CRP row 34
$ws.Range("H34").Value = 18813.83
$ws.Range("I34").Value = 46584.227
$ws.Range("J34").Value = 2301.7026
$ws.Range("K34").Value = 46584.227
$ws.Range("L34").Value = 2301.7026
$ws.Range("M34").Value = -46382.227
$ws.Range("N34").Value = -2705.7026

# CRP row 133
$ws.Range("H133").Value = 63470
$ws.Range("J133").Value = 63470
$ws.Range("L133").Value = 63470
$ws.Range("N133").Value = -68530

# CRP row 134
$ws.Range("H134").Value = 1453.65
$ws.Range("I134").Value = 1121.1177
$ws.Range("J134").Value = 3338
$ws.Range("K134").Value = 3363.3531
$ws.Range("L134").Value = 10014
$ws.Range("M134").Value = -828.3531000000003
$ws.Range("N134").Value = -15084

# CRP row 141
$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# CUL row 17
$ws.Range("H17").Value = 2161.111
$ws.Range("I17").Value = 2214.2856
$ws.Range("K17").Value = 6642.8568
$ws.Range("M17").Value = -6473.8568

# CUL row 23
$ws.Range("H23").Value = 804.2
$ws.Range("I23").Value = 150
$ws.Range("J23").Value = 876.8889
$ws.Range("K23").Value = 450
$ws.Range("L23").Value = 2630.6667
$ws.Range("M23").Value = -215
$ws.Range("N23").Value = -3100.6667

# CUL row 38
$ws.Range("H38").Value = 93.85714
$ws.Range("I38").Value = 103.666664
$ws.Range("J38").Value = 86.5
$ws.Range("K38").Value = 310.999992
$ws.Range("L38").Value = 259.5
$ws.Range("M38").Value = 36.00000799999998
$ws.Range("N38").Value = -953.5

# CUL row 44
$ws.Range("H44").Value = 699.6667
$ws.Range("I44").Value = 549.5
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 1648.5
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = -1250.5
$ws.Range("N44").Value = -3796

# CUL row 134
$ws.Range("H134").Value = 3968.3333
$ws.Range("I134").Value = 2526.6667
$ws.Range("K134").Value = 7580.000100000001
$ws.Range("M134").Value = -2510.000100000001

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Range("H2").Value = 264.8421
$ws.Range("I2").Value = 247.86667
$ws.Range("J2").Value = 328.5
$ws.Range("K2").Value = 247.86667
$ws.Range("L2").Value = 328.5
$ws.Range("M2").Value = -134.86667
$ws.Range("N2").Value = -554.5

# GSM row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

# GSM row 62
$ws.Range("H62").Value = 14519.25
$ws.Range("I62").Value = 13077
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 13077
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -12391
$ws.Range("N62").Value = -16372

# GSM row 65
$ws.Range("H65").Value = 14519.25
$ws.Range("I65").Value = 13077
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 39231
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -35799
$ws.Range("N65").Value = -51864

# GSM row 102
$ws.Range("H102").Value = 275230.8
$ws.Range("I102").Value = 3058.2144
$ws.Range("K102").Value = 3058.2144
$ws.Range("M102").Value = -1436.2144

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 481.46667
$ws.Range("I22").Value = 426.66666
$ws.Range("K22").Value = 426.66666
$ws.Range("M22").Value = -131.66666

# LTW row 27
$ws.Range("H27").Value = 481.46667
$ws.Range("I27").Value = 426.66666
$ws.Range("K27").Value = 426.66666
$ws.Range("M27").Value = -319.66666

# LTW row 55
$ws.Range("H55").Value = 494500.53
$ws.Range("I55").Value = 811859.5600000001
$ws.Range("J55").Value = 830.8889
$ws.Range("K55").Value = 811859.5600000001
$ws.Range("L55").Value = 830.8889
$ws.Range("M55").Value = -811686.5600000001
$ws.Range("N55").Value = -1176.8889

# LTW row 96
$ws.Range("H96").Value = 15747.667
$ws.Range("J96").Value = 15747.667
$ws.Range("L96").Value = 15747.667
$ws.Range("N96").Value = -21239.667

$ws = $wb.Worksheets.Item("WVR")
# WVR row 38
$ws.Range("H38").Value = 6144.25
$ws.Range("I38").Value = 3232
$ws.Range("J38").Value = 6560.2856
$ws.Range("K38").Value = 3232
$ws.Range("L38").Value = 6560.2856
$ws.Range("M38").Value = -2759
$ws.Range("N38").Value = -7506.2856

# WVR row 132
$ws.Range("H132").Value = 3097.575
$ws.Range("I132").Value = 3004.2646
$ws.Range("J132").Value = 3626.3333
$ws.Range("K132").Value = 9012.793799999999
$ws.Range("L132").Value = 10878.9999
$ws.Range("M132").Value = -6482.793799999999
$ws.Range("N132").Value = -15938.9999

# WVR row 139
$ws.Range("H139").Value = 65657.5
$ws.Range("J139").Value = 65657.5
$ws.Range("L139").Value = 65657.5
$ws.Range("N139").Value = -75937.5
